$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 63

$ws.Range("C3").Value = 25
$ws.Range("D3").Value = 0.25
$ws.Range("E3").Value = 0.2
$ws.Range("F3").Value = 0.2222222222222222

$ws.Range("C4").Value = 26

$ws.Range("C5").Value = 30
$ws.Range("D5").Value = 0.3548387096774194
$ws.Range("E5").Value = 0.3666666666666666
$ws.Range("F5").Value = 0.360655737704918

$ws.Range("C6").Value = 56
$ws.Range("D6").Value = 0.2697368421052632
$ws.Range("E6").Value = 0.7321428571428571
$ws.Range("F6").Value = 0.3942307692307692

$ws.Range("C7").Value = 9
